# Fruta / hortaliza, semanal
#
# The source table (Cilantro @ Mercado Mayorista Lo Valledor de Santiago) gets
# two additional weekly records inserted right after the existing row for
# 2022-01-08/15ish (row 650), pushing every subsequent record down by two rows
# and appending the same two records that used to trail off the bottom of the
# sheet (so the sheet grows from 691 to 693 data-ish rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the current row 651, shifting rows 651:691
# down to 653:693 (and pushing the dimension to A1:R693).
$ws.Rows.Item(651).Resize(2).Insert()

# Row 651 — new weekly record ("Primera", $/caja 36 atados)
$ws.Range("A651").Value = 6
$ws.Range("B651").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C651").Value = "Metropolitana"
$ws.Range("D651").Value = 44585
$ws.Range("E651").Value = 13
$ws.Range("F651").Value = 100112040
$ws.Range("G651").Value = "Cilantro"
$ws.Range("H651").Value = "Sin especificar"
$ws.Range("I651").Value = "Primera"
$ws.Range("J651").Value = 450
$ws.Range("K651").Value = 7500
$ws.Range("L651").Value = 8000
$ws.Range("M651").Value = 7711
$ws.Range("N651").Value = "$/caja 36 atados"
$ws.Range("O651").Value = "Región Metropolitana"
$ws.Range("P651").Value = 214
$ws.Range("Q651").Value = 36
$ws.Range("R651").Value = "Hortaliza"

# Row 652 — new weekly record ("Primera", $/docena de atados)
$ws.Range("A652").Value = 6
$ws.Range("B652").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C652").Value = "Metropolitana"
$ws.Range("D652").Value = 44585
$ws.Range("E652").Value = 13
$ws.Range("F652").Value = 100112040
$ws.Range("G652").Value = "Cilantro"
$ws.Range("H652").Value = "Sin especificar"
$ws.Range("I652").Value = "Primera"
$ws.Range("J652").Value = 320
$ws.Range("K652").Value = 15000
$ws.Range("L652").Value = 16000
$ws.Range("M652").Value = 15406
$ws.Range("N652").Value = "$/docena de atados"
$ws.Range("O652").Value = "Región Metropolitana"
$ws.Range("P652").Value = 5135
$ws.Range("Q652").Value = 3
$ws.Range("R652").Value = "Hortaliza"
